$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (fill/style) from existing alternating-color player blocks
# so the new rows reuse the same style indexes (s="2" / s="3") as the
# pre-existing data, matching the green/yellow banding pattern.

# Rows 32-34: Christian Rozeboom -> same style as Chris Board (rows 2-4, s="2")
$ws.Range("A2:F4").Copy()
$ws.Range("A32:F34").PasteSpecial(-4122)

# Rows 35-37: Duke Riley -> same style as Brian Burns (rows 5-7, s="3")
$ws.Range("A5:F7").Copy()
$ws.Range("A35:F37").PasteSpecial(-4122)

# Rows 38-40: Troy Reeder -> same style as Chris Board (rows 2-4, s="2")
$ws.Range("A2:F4").Copy()
$ws.Range("A38:F40").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 32: Christian Rozeboom / Group1
$ws.Range("A32").Value = "Christian Rozeboom"
$ws.Range("B32").Value = "Group1"
$ws.Range("C32").Value = 2.666666666666667
$ws.Range("D32").Value = 76.66666666666667
$ws.Range("E32").Value = 45.66666666666666
$ws.Range("F32").Value = 31

# Row 33: Christian Rozeboom / Group2
$ws.Range("A33").Value = "Christian Rozeboom"
$ws.Range("B33").Value = "Group2"
$ws.Range("C33").Value = 0.6666666666666666
$ws.Range("D33").Value = 26.66666666666667
$ws.Range("E33").Value = 14.33333333333333
$ws.Range("F33").Value = 12.33333333333333

# Row 34: Christian Rozeboom / Difference
$ws.Range("A34").Value = "Christian Rozeboom"
$ws.Range("B34").Value = "Difference"
$ws.Range("C34").Value = -2
$ws.Range("D34").Value = -50
$ws.Range("E34").Value = -31.33333333333333
$ws.Range("F34").Value = -18.66666666666666

# Row 35: Duke Riley / Group1
$ws.Range("A35").Value = "Duke Riley"
$ws.Range("B35").Value = "Group1"
$ws.Range("C35").Value = 0.3333333333333333
$ws.Range("D35").Value = 30.11111111111111
$ws.Range("E35").Value = 17.66666666666667
$ws.Range("F35").Value = 12.44444444444444

# Row 36: Duke Riley / Group2
$ws.Range("A36").Value = "Duke Riley"
$ws.Range("B36").Value = "Group2"
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 32.66666666666666
$ws.Range("E36").Value = 20
$ws.Range("F36").Value = 12.66666666666667

# Row 37: Duke Riley / Difference
$ws.Range("A37").Value = "Duke Riley"
$ws.Range("B37").Value = "Difference"
$ws.Range("C37").Value = 0.6666666666666667
$ws.Range("D37").Value = 2.555555555555554
$ws.Range("E37").Value = 2.333333333333332
$ws.Range("F37").Value = 0.2222222222222232

# Row 38: Troy Reeder / Group1
$ws.Range("A38").Value = "Troy Reeder"
$ws.Range("B38").Value = "Group1"
$ws.Range("C38").Value = 2.666666666666667
$ws.Range("D38").Value = 76.66666666666667
$ws.Range("E38").Value = 45.66666666666666
$ws.Range("F38").Value = 31

# Row 39: Troy Reeder / Group2
$ws.Range("A39").Value = "Troy Reeder"
$ws.Range("B39").Value = "Group2"
$ws.Range("C39").Value = 0.6666666666666666
$ws.Range("D39").Value = 26.66666666666667
$ws.Range("E39").Value = 14.33333333333333
$ws.Range("F39").Value = 12.33333333333333

# Row 40: Troy Reeder / Difference
$ws.Range("A40").Value = "Troy Reeder"
$ws.Range("B40").Value = "Difference"
$ws.Range("C40").Value = -2
$ws.Range("D40").Value = -50
$ws.Range("E40").Value = -31.33333333333333
$ws.Range("F40").Value = -18.66666666666666
